# Added a new US state chart excluding New York to zoom in on the other states.
#
# I14 becomes a hard input value (actual reported US case count for that day)
# instead of a projected formula, and the forward-projection formulas in
# I15:I28 widen their trailing AVERAGE() window from 3 days to 5 days
# (M(n-5):M(n-1) instead of M(n-3):M(n-1)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I14: replace the formula with the literal input value, and restyle it like
# the other manually-entered "I" cells above it (I10:I13).
$ws.Range("I10").Copy()
$ws.Range("I14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I14").Value = 19383

# I15: widen the trailing average window to 5 days.
$ws.Range("I15").Formula = "=I14*(1+AVERAGE(M10:M14))"

# I16:I28: same widened-window pattern, written to the whole block at once so
# Excel groups it as one shared formula (relative refs shift per row, same as
# dragging the fill handle down from I16).
$ws.Range("I16:I28").Formula = "=I15*(1+AVERAGE(M11:M15))"

# Mark the active cell/selection as it was left after the edit.
$ws.Range("I26").Select()
